$wb = $excel.ActiveWorkbook
$sourceSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sourceSheet.Copy($null, $sourceSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2025-08-28"

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "魔術師クノンは見えている"
$ws.Cells.Item(2, 3).Value = "La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)"
$ws.Cells.Item(2, 4).Value = "第40話①"

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "実は俺、最強でした？"
$ws.Cells.Item(3, 3).Value = "原作：澄守 彩 漫画：高橋 愛"
$ws.Cells.Item(3, 4).Value = "第123話　王妃とハルト・後編"

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "生徒会にも穴はある！"
$ws.Cells.Item(4, 3).Value = "むちまろ"
$ws.Cells.Item(4, 4).Value = "第134話`t太賀のトラウマ"

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "男女比1：5の世界でも普通に生きられると思った？　～激重感情な彼女たちが無自覚男子に翻弄されたら～"
$ws.Cells.Item(5, 3).Value = "三藤 孝太郎(原作) 桃季憂(漫画) jimmy(キャラクター原案)"
$ws.Cells.Item(5, 4).Value = "第10話-2"

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "ダンジョンの幼なじみ"
$ws.Cells.Item(6, 3).Value = "久真やすひさ(著者)"
$ws.Cells.Item(6, 4).Value = "第56話"

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "異世界魔王と召喚少女の奴隷魔術"
$ws.Cells.Item(7, 3).Value = "原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大"
$ws.Cells.Item(7, 4).Value = "『異世界魔王』ヒロイン総選挙 結果発表"

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～"
$ws.Cells.Item(8, 3).Value = "漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき"
$ws.Cells.Item(8, 4).Value = "第33話 独身貴族は見積もりを誤る（2）"

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "北斗の拳 世紀末ドラマ撮影伝"
$ws.Cells.Item(9, 3).Value = "原案/武論尊・原哲夫 漫画/倉尾宏"
$ws.Cells.Item(9, 4).Value = "第75話 宙を舞う悪役俳優‼︎"

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "スキル【万物支配】に目覚めたおっさんは、ダンジョンで生計を立てることにしました～無職から始める支配者無双～"
$ws.Cells.Item(10, 3).Value = "岸本和葉 原田 臙 シミズヒロノリ 吉武"
$ws.Cells.Item(10, 4).Value = "第5話　一撃(後編)"

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "ダンジョン・シェルパ　迷宮道先案内人"
$ws.Cells.Item(11, 3).Value = "原作/加茂セイ 漫画/刀坂アキラ"
$ws.Cells.Item(11, 4).Value = "潜行：61(前編)"

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "オタクに優しいギャルはいない!?"
$ws.Cells.Item(12, 3).Value = "のりしろちゃん 魚住さかな"
$ws.Cells.Item(12, 4).Value = "【#154】起きたーーッ"

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "生徒会役員共"
$ws.Cells.Item(13, 3).Value = "氏家ト全"
$ws.Cells.Item(13, 4).Value = "#411"

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "カナン様はあくまでチョロい"
$ws.Cells.Item(14, 3).Value = "nonco"
$ws.Cells.Item(14, 4).Value = "第148話`tカナンの布団の中"

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "脇役に転生したはずが、いつの間にか伝説の錬金術師になってた～仲間たちが英雄でも俺は支援職なんだが～"
$ws.Cells.Item(15, 3).Value = "神無月みり 相野 仁"
$ws.Cells.Item(15, 4).Value = "第２７話　脇役、先輩の無事を祈る（１）"

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "29歳独身は異世界で自由に生きた……かった。"
$ws.Cells.Item(16, 3).Value = "オオハマイコ(漫画) リュート(原作) 桑島黎音(キャラクター原案)"
$ws.Cells.Item(16, 4).Value = "第41話-2②"

$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "田舎のホームセンター男の自由な異世界生活"
$ws.Cells.Item(17, 3).Value = "うさぴょん(原作) 古来歩(漫画) 市丸きすけ(キャラクター原案)"
$ws.Cells.Item(17, 4).Value = "第76話その1"

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "ポンコツ風紀委員とスカート丈が不適切なＪＫの話"
$ws.Cells.Item(18, 3).Value = "横田卓馬"
$ws.Cells.Item(18, 4).Value = "番外編　ポンコツ一年生たちの文化祭準備の話"

$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "まんきつしたい常連さん"
$ws.Cells.Item(19, 3).Value = "しんみりん(著者)"
$ws.Cells.Item(19, 4).Value = "第47話前編"

$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "すべての人類を破壊する。それらは再生できない。"
$ws.Cells.Item(20, 3).Value = "横田卓馬(漫画) 伊瀬勝良(原作)"
$ws.Cells.Item(20, 4).Value = "第68話その1"

$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。"
$ws.Cells.Item(21, 3).Value = "マツモトケンゴ"
$ws.Cells.Item(21, 4).Value = "第６３話　ダンスゲームの戦いが始まった（２）"

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―"
$ws.Cells.Item(22, 3).Value = "光永康則"
$ws.Cells.Item(22, 4).Value = "第６８話『施錠停止』④"

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "えろいことするために巨乳美少女奴隷を買ったはずが、お師匠さまと慕われて思った通りにいかなくなる話"
$ws.Cells.Item(23, 3).Value = "佐藤36(作画) 煮豆シューター(原作)"
$ws.Cells.Item(23, 4).Value = "第6話前半"

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "二番目な僕と一番の彼女"
$ws.Cells.Item(24, 3).Value = "ぬずタニ(漫画) 和尚(原作) ミュシャ(キャラクター原案)"
$ws.Cells.Item(24, 4).Value = "第2話②"

$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "転生錬金少女のスローライフ"
$ws.Cells.Item(25, 3).Value = "里町(漫画) 夜想庭園(原作) potg(キャラクター原案)"
$ws.Cells.Item(25, 4).Value = "第17話-3"

$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "善人おっさん、生まれ変わったらSSSランク人生が確定した"
$ws.Cells.Item(26, 3).Value = "原作／三木なずな 漫画／ゆづましろ キャラクター原案／伍長"
$ws.Cells.Item(26, 4).Value = "第134話"

$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "よわよわ先生"
$ws.Cells.Item(27, 3).Value = "福地カミオ"
$ws.Cells.Item(27, 4).Value = "第113話`tエモエモの最終課題③"

$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "ギャルゲーマーに褒められたい"
$ws.Cells.Item(28, 3).Value = "げしゅまろ(著者)"
$ws.Cells.Item(28, 4).Value = "47話"

$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "転移したら奴隷の父になったけど、家族として愛でることにした〜実は勇者と魔王と聖女だった子供たちの力がとんでもスキルで使えたので最強です〜"
$ws.Cells.Item(29, 3).Value = "えむだ(作画) 御峰。(原作)"
$ws.Cells.Item(29, 4).Value = "第3話前半"

$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "百瀬アキラの初恋破綻中。"
$ws.Cells.Item(30, 3).Value = "晴川シンタ"
$ws.Cells.Item(30, 4).Value = "第37話 真夏の海辺に腰かけ中。"

$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "無能は不要と言われ『時計使い』の僕は職人ギルドから追い出されるも、ダンジョンの深部で真の力に覚醒する"
$ws.Cells.Item(31, 3).Value = "漫画：さらさみさ 小説： 桜霧琥珀 キャラ原案： 福きつね"
$ws.Cells.Item(31, 4).Value = "第19話後半"

$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "勇者パーティーの荷物持ち"
$ws.Cells.Item(32, 3).Value = "原作：河本ほむら／作画：八嶋諒"
$ws.Cells.Item(32, 4).Value = "第25話 荷物持ちと勇者パーティーの戦士②"

$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜"
$ws.Cells.Item(33, 3).Value = "戸賀 環 坂木持丸 riritto"
$ws.Cells.Item(33, 4).Value = "第51話②　呪われた家を探索してみた"

$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "幼女戦記"
$ws.Cells.Item(34, 3).Value = "東條チカ(漫画) カルロ・ゼン(原作) 篠月しのぶ(キャラクター原案)"
$ws.Cells.Item(34, 4).Value = "第百七章：ドードーバード航空戦Ⅱ"

$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "物語の黒幕に転生して"
$ws.Cells.Item(35, 3).Value = "瀬川はじめ(漫画) 結城涼(原作) なかむら(キャラクター原案)"
$ws.Cells.Item(35, 4).Value = "第33話"

$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "Lv２からチートだった元勇者候補のまったり異世界ライフ"
$ws.Cells.Item(36, 3).Value = "糸町秋音（漫画） 鬼ノ城ミヤ（原作） 片桐（キャラクター原案）"
$ws.Cells.Item(36, 4).Value = "第60話　ふたつの希望（後編）"

$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "シャドウ・アサシンズ・ワールド ～影は薄いけど、最強忍者やってます～"
$ws.Cells.Item(37, 3).Value = "空山トキ 五色安未 泉乃せん 伍長"
$ws.Cells.Item(37, 4).Value = "第12話　少女と本当の自分３（1）"

$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "ありふれた職業で世界最強"
$ws.Cells.Item(38, 3).Value = "RoGa（漫画） 白米 良（原作） たかやKi（キャラクター原案）"
$ws.Cells.Item(38, 4).Value = "第84話　人間らしさ（後編）"

$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "おかけになった呪文は現在使われておりません"
$ws.Cells.Item(39, 3).Value = "ロケット商会 天宮ケイリ"
$ws.Cells.Item(39, 4).Value = "第4話　”あの男”"

$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "絶対死なないステラ姫"
$ws.Cells.Item(40, 3).Value = "光永康則 大高稲"
$ws.Cells.Item(40, 4).Value = "第１５話　絶対指名手配されない（１）"

$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "葉木莉さんは君だけの死神になりたい"
$ws.Cells.Item(41, 3).Value = "35まち"
$ws.Cells.Item(41, 4).Value = "第７話"

$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "ひとりぼっちの異世界攻略"
$ws.Cells.Item(42, 3).Value = "びび（漫画） 五示正司（原作）"
$ws.Cells.Item(42, 4).Value = "第233話　ミラクルな幕引き"

$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "DT転生　～３０歳まで童貞で転生したら、史上最強の魔法使いになりました！～"
$ws.Cells.Item(43, 3).Value = "石田衣良 山田秋太郎"
$ws.Cells.Item(43, 4).Value = "第２１話　合言葉は「アーチボルト」！（２）"

$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "調教師は魔物に囲まれて生きていきます。～勇者パーティーに置いていかれたけど、伝説の魔物と出会い最強になってた～"
$ws.Cells.Item(44, 3).Value = "尾切美月(作画) 七篠 龍(原作) 竹花ノート(キャラクター原案)"
$ws.Cells.Item(44, 4).Value = "第6話前半"

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "貴族令嬢がジャンクフード食って「美味いですわ！」するだけの話"
$ws.Cells.Item(45, 3).Value = "ごくげつ(作画) パイルバンカー串山(原作)"
$ws.Cells.Item(45, 4).Value = "第6話前半"

$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "おはらい箱の天才付与術師は、辺境で悠々自適に暮らしたい～万能付与術で気付いたら辺境が世界最強の快適拠点になっていた～"
$ws.Cells.Item(46, 3).Value = "柊木 楸(作画) 水無月(原作) 布施龍太(キャラクター原案)"
$ws.Cells.Item(46, 4).Value = "第6話前半"

$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "アンゴルモア 元寇合戦記　【博多編】"
$ws.Cells.Item(47, 3).Value = "たかぎ七彦(著者)"
$ws.Cells.Item(47, 4).Value = "第四十五話その四"

$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する"
$ws.Cells.Item(48, 3).Value = "作画：マエD 原作：新人"
$ws.Cells.Item(48, 4).Value = "第6話(2)"

$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "黒月のイェルクナハト"
$ws.Cells.Item(49, 3).Value = "スズモトコウ"
$ws.Cells.Item(49, 4).Value = "第8話`t死にに行け"

$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "帰ってください！ 阿久津さん"
$ws.Cells.Item(50, 3).Value = "長岡太一(著者)"
$ws.Cells.Item(50, 4).Value = "第195話"

$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信"
$ws.Cells.Item(51, 3).Value = "下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)"
$ws.Cells.Item(51, 4).Value = "第2話前編"
